# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos sheet
# with the latest scraped values (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.242.17'
$ws.Range('E2').Value = '  +7.55%  '
$ws.Range('D3').Value = '3.034.22'
$ws.Range('E3').Value = '  +5.29%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '''582.80'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.86%  '
$ws.Range('D6').Value = '''157.97'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +11.16%  '
$ws.Range('D8').Value = '3.028.41'
$ws.Range('E8').Value = '  +5.15%  '
$ws.Range('E9').Value = '  +3.72%  '
$ws.Range('E10').Value = '  +2.11%  '
$ws.Range('E11').Value = '  +6.96%  '
$ws.Range('D12').Value = '''0.453'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.77%  '
$ws.Range('D13').Value = '''0.0000252'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +9.88%  '
$ws.Range('D14').Value = '''34.64'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +9.67%  '
$ws.Range('E15').Value = '  +0.79%  '
$ws.Range('D16').Value = '66.334.65'
$ws.Range('E16').Value = '  +7.76%  '
$ws.Range('D17').Value = '3.541.02'
$ws.Range('E17').Value = '  +5.40%  '
$ws.Range('D18').Value = '''6.96'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +7.32%  '
$ws.Range('D19').Value = '3.032.90'
$ws.Range('E19').Value = '  +4.73%  '
$ws.Range('D20').Value = '''464.28'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +8.23%  '
$ws.Range('D21').Value = '''13.89'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.02%  '
$ws.Range('E22').Value = '  +5.72%  '
$ws.Range('D23').Value = '''7.37'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.62%  '
$ws.Range('D24').Value = '''82.53'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.78%  '
$ws.Range('D25').Value = '''2.25'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +12.46%  '
$ws.Range('D26').Value = '''12.48'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.37%  '
$ws.Range('D27').Value = '''10.63'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.62%  '
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('D29').Value = '''8.03'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +14.34%  '
$ws.Range('D30').Value = '''2.35'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +16.74%  '
$ws.Range('D31').Value = '''0.0000105'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('E32').Value = '  +5.00%  '
$ws.Range('D33').Value = '''27.10'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.95%  '
$ws.Range('E34').Value = '  +5.54%  '
$ws.Range('D35').Value = '''0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('D36').Value = '''0.999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.70%  '
$ws.Range('D37').Value = '''5.79'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.19%  '
$ws.Range('E38').Value = '  +14.92%  '
$ws.Range('D39').Value = '''3.06'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +9.85%  '
$ws.Range('D40').Value = '''49.61'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.73%  '
$ws.Range('E41').Value = '  +7.82%  '
$ws.Range('D42').Value = '''43.89'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +11.98%  '
$ws.Range('D43').Value = '''0.302'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +14.11%  '
$ws.Range('D44').Value = '''8.45'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.92%  '
$ws.Range('D45').Value = '''390.37'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +13.51%  '
$ws.Range('D46').Value = '2.812.42'
$ws.Range('E46').Value = '  +5.04%  '
$ws.Range('D47').Value = '''0.0354'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +6.20%  '
$ws.Range('D48').Value = '''134.17'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.66%  '
$ws.Range('D50').Value = '''23.56'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +10.40%  '
